# Apply the "adding term 2.0 now utf-8" update:
#  - bump Version / Date / Contact on the Metadata sheet
#  - add a second copy of the "Include from FSIII" sheet, named
#    "Include from FSIII 2", at the end of the workbook
#  - the existing "Include from FSIII" sheet's Operation/B row gets its
#    Value cell turned into the new term's tracking GUID (matches the
#    source edit exactly)

$wb = $excel.ActiveWorkbook

# --- 1. Update the Metadata sheet -----------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- 2. Duplicate "Include from FSIII" to the end of the workbook ---------
$inc = $wb.Worksheets.Item("Include from FSIII")
$inc.Copy($null, $inc)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Include from FSIII 2"

# --- 3. Update the original sheet's Value cell for the "B" concept row ----
$inc.Range("C2").Value = "2c02a704-deee-4878-9378-1167613b3da6"

# --- 4. Keep the Metadata sheet as the active/selected tab -----------------
$meta.Activate()
